# edit.ps1 -- apply the "schedule dates" commit to Schdule.docx
#
# Summary of the change (see commit message / diff):
#   - Append " – March 9"  to the "Modeling of the maps" bullet (Physics section)
#   - Append " – March 30" to the "Power Ups" bullet (Physics section)
#   - Append " – Feb 3"    to the "Spring Implementation" bullet
#   - Append " – Feb 3"    to the "Collisions" bullet
#   - Append " - Ongoing"  to the "Integration with graphics / rendering" bullet
#   - Remove the "Ability to switch between simulation and faking physics" bullet
#     entirely (it is merged away / deleted in the diff).
#
# (The diff also shows Word's automatic spell/grammar-checker proofErr markers
#  being added around "shaders", "MultiView" and "Displaying" -- those are
#  cosmetic proofing artifacts with no text-content change, so there is
#  nothing to edit there.)

$d = $word.ActiveDocument

function Get-ParagraphByText($needle) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $t = $p.Range.Text.TrimEnd([char]13).TrimEnd(" ")
        # Use .Equals() (ordinal/case-sensitive) rather than -eq, which is
        # case-insensitive and would conflate e.g. "Power Ups" / "Power ups".
        if ($t.Equals($needle)) {
            return $p
        }
    }
    return $null
}

function Append-TextToParagraph($needle, $suffix) {
    $p = Get-ParagraphByText $needle
    if ($p -eq $null) {
        Write-Output "WARNING: paragraph not found: $needle"
        return
    }
    $endPos = $p.Range.End - 1   # just before the paragraph mark
    $insPt = $d.Range($endPos, $endPos)
    $insPt.InsertBefore($suffix)
}

# 1. "Modeling of the maps" -> append " – March 9"
Append-TextToParagraph "Modeling of the maps" " $([char]0x2013) March 9"

# 2. "Power Ups" -> append " – March 30"
Append-TextToParagraph "Power Ups" " $([char]0x2013) March 30"

# 3. "Spring Implementation" -> append " – Feb 3"
Append-TextToParagraph "Spring Implementation" " $([char]0x2013) Feb 3"

# 4. "Collisions" -> append " – Feb 3"
Append-TextToParagraph "Collisions" " $([char]0x2013) Feb 3"

# 5. "Integration with graphics / rendering" -> append " - Ongoing"
Append-TextToParagraph "Integration with graphics / rendering" " - Ongoing"

# 6. Remove the "Ability to switch between simulation and faking physics" bullet entirely.
$p = Get-ParagraphByText "Ability to switch between simulation and faking physics"
if ($p -ne $null) {
    $p.Range.Delete()
} else {
    Write-Output "WARNING: paragraph not found: Ability to switch between simulation and faking physics"
}
